# Scheduled-runner update: refresh cached market-price / profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) across the per-job Leve
# profit sheets, per the latest price pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("N18").ClearContents()
$ws.Range("H99").Value = 20163.5
$ws.Range("I99").Value = 20163.5
$ws.Range("K99").Value = 60490.5
$ws.Range("M99").Value = -58992.5
$ws.Range("H138").Value = 2855.242
$ws.Range("J138").Value = 4889.615
$ws.Range("L138").Value = 14668.845
$ws.Range("N138").Value = -24948.845

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13335578
$ws.Range("I32").Value = 13890700
$ws.Range("K32").Value = 13890700
$ws.Range("M32").Value = -13890413
$ws.Range("H63").Value = 1998.7142
$ws.Range("I63").Value = 1830.8334
$ws.Range("K63").Value = 1830.8334
$ws.Range("M63").Value = -1144.8334
$ws.Range("H66").Value = 1998.7142
$ws.Range("I66").Value = 1830.8334
$ws.Range("K66").Value = 9154.166999999999
$ws.Range("M66").Value = -5722.166999999999
$ws.Range("H97").Value = 1078.6875
$ws.Range("I97").Value = 1083.9333
$ws.Range("K97").Value = 1083.9333
$ws.Range("M97").Value = -587.9332999999999
$ws.Range("H122").Value = 38464092
$ws.Range("I122").Value = 2709.4443
$ws.Range("K122").Value = 8128.3329
$ws.Range("M122").Value = -5678.3329
$ws.Range("H139").Value = 68749.25
$ws.Range("J139").Value = 68749.25
$ws.Range("L139").Value = 68749.25
$ws.Range("N139").Value = -79029.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 863.3333
$ws.Range("I22").Value = 795
$ws.Range("K22").Value = 795
$ws.Range("M22").Value = -622
$ws.Range("H94").Value = 4125
$ws.Range("I94").Value = 3250
$ws.Range("J94").Value = 5000
$ws.Range("K94").Value = 3250
$ws.Range("L94").Value = 5000
$ws.Range("M94").Value = -2799
$ws.Range("N94").Value = -5902
$ws.Range("H134").Value = 2012.1
$ws.Range("I134").Value = 976.04
$ws.Range("K134").Value = 2928.12
$ws.Range("M134").Value = -393.1199999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 143.5
$ws.Range("I5").Value = 143.5
$ws.Range("K5").Value = 143.5
$ws.Range("M5").Value = -31.5
$ws.Range("H22").Value = 1029.7273
$ws.Range("I22").Value = 384.8
$ws.Range("J22").Value = 1567.1666
$ws.Range("K22").Value = 384.8
$ws.Range("L22").Value = 1567.1666
$ws.Range("M22").Value = -34.80000000000001
$ws.Range("N22").Value = -2267.1666
$ws.Range("H31").Value = 35162.574
$ws.Range("I31").Value = 1453
$ws.Range("K31").Value = 1453
$ws.Range("M31").Value = -1158
$ws.Range("H34").Value = 35162.574
$ws.Range("I34").Value = 1453
$ws.Range("K34").Value = 1453
$ws.Range("M34").Value = -1251
$ws.Range("H58").Value = 3499.04
$ws.Range("I58").Value = 1437.6471
$ws.Range("K58").Value = 1437.6471
$ws.Range("M58").Value = -1234.6471
$ws.Range("H99").Value = 4687.4375
$ws.Range("I99").Value = 4727.1816
$ws.Range("J99").Value = 4600
$ws.Range("K99").Value = 4727.1816
$ws.Range("L99").Value = 4600
$ws.Range("M99").Value = -3229.1816
$ws.Range("N99").Value = -7596
$ws.Range("H126").Value = 4687.4375
$ws.Range("I126").Value = 4727.1816
$ws.Range("J126").Value = 4600
$ws.Range("K126").Value = 14181.5448
$ws.Range("L126").Value = 13800
$ws.Range("M126").Value = -11711.5448
$ws.Range("N126").Value = -18740
$ws.Range("H131").Value = 46799.8
$ws.Range("J131").Value = 46249.75
$ws.Range("L131").Value = 46249.75
$ws.Range("N131").Value = -56329.75
$ws.Range("H132").Value = 2101
$ws.Range("I132").Value = 1776
$ws.Range("K132").Value = 5328
$ws.Range("M132").Value = -2798
$ws.Range("H134").Value = 5578.375
$ws.Range("I134").Value = 3515.75
$ws.Range("J134").Value = 9016.083000000001
$ws.Range("K134").Value = 10547.25
$ws.Range("L134").Value = 27048.249
$ws.Range("M134").Value = -8012.25
$ws.Range("N134").Value = -32118.249
$ws.Range("H136").Value = 3499.04
$ws.Range("I136").Value = 1437.6471
$ws.Range("K136").Value = 4312.9413
$ws.Range("M136").Value = -1762.9413
$ws.Range("H140").Value = 88671.25
$ws.Range("J140").Value = 88671.25
$ws.Range("L140").Value = 88671.25
$ws.Range("N140").Value = -99031.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2762.4375
$ws.Range("J5").Value = 3593.3
$ws.Range("L5").Value = 10779.9
$ws.Range("N5").Value = -11003.9
$ws.Range("H32").Value = 9901.416999999999
$ws.Range("J32").Value = 11601.7
$ws.Range("L32").Value = 34805.10000000001
$ws.Range("N32").Value = -35371.10000000001
$ws.Range("H70").Value = 11623
$ws.Range("H73").Value = 11623
$ws.Range("H131").Value = 6390662.5
$ws.Range("J131").Value = 5118811.5
$ws.Range("L131").Value = 15356434.5
$ws.Range("N131").Value = -15366514.5
$ws.Range("H134").Value = 15878860
$ws.Range("I134").Value = 1048.5
$ws.Range("J134").Value = 30313234
$ws.Range("K134").Value = 3145.5
$ws.Range("L134").Value = 90939702
$ws.Range("M134").Value = 1924.5
$ws.Range("N134").Value = -90949842
$ws.Range("H135").Value = 2762.4375
$ws.Range("J135").Value = 3593.3
$ws.Range("L135").Value = 32339.7
$ws.Range("N135").Value = -37409.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 36999
$ws.Range("J53").Value = 36999
$ws.Range("L53").Value = 36999
$ws.Range("N53").Value = -38261
$ws.Range("H132").Value = 440894.66
$ws.Range("I132").Value = 455960
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 1367880
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -1365350
$ws.Range("N132").Value = -17060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 9139.154
$ws.Range("I40").Value = 9644.888999999999
$ws.Range("K40").Value = 9644.888999999999
$ws.Range("M40").Value = -9508.888999999999
$ws.Range("H46").Value = 3736.2727
$ws.Range("I46").Value = 832.6667
$ws.Range("J46").Value = 4825.125
$ws.Range("K46").Value = 832.6667
$ws.Range("L46").Value = 4825.125
$ws.Range("M46").Value = -644.6667
$ws.Range("N46").Value = -5201.125
$ws.Range("H136").Value = 6017.3335
$ws.Range("I136").Value = 2946.3076
$ws.Range("J136").Value = 14002
$ws.Range("K136").Value = 8838.9228
$ws.Range("L136").Value = 42006
$ws.Range("M136").Value = -6288.9228
$ws.Range("N136").Value = -47106

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1021.64703
$ws.Range("I107").Value = 757.3
$ws.Range("J107").Value = 1399.2858
$ws.Range("K107").Value = 2271.9
$ws.Range("L107").Value = 4197.857400000001
$ws.Range("M107").Value = -351.8999999999996
$ws.Range("N107").Value = -8037.857400000001
$ws.Range("H122").Value = 4356
$ws.Range("I122").Value = 3229.1538
$ws.Range("K122").Value = 9687.4614
$ws.Range("M122").Value = -7237.4614
$ws.Range("H132").Value = 4708.254
$ws.Range("I132").Value = 2195.9788
$ws.Range("J132").Value = 12088.0625
$ws.Range("K132").Value = 6587.9364
$ws.Range("L132").Value = 36264.1875
$ws.Range("M132").Value = -4057.9364
$ws.Range("N132").Value = -41324.1875
$ws.Range("H136").Value = 2643.4666
$ws.Range("I136").Value = 1874.174
$ws.Range("J136").Value = 5171.143
$ws.Range("K136").Value = 5622.522
$ws.Range("L136").Value = 15513.429
$ws.Range("M136").Value = -3072.522
$ws.Range("N136").Value = -20613.429
